# "moved excel to pages and removed expected outcome"
# Net content change observed in the target workbook: the "python DS" sheet
# gains a second column ("expected Outcome") describing what each code
# snippet in column A should produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python DS")
$ws.Activate()

# New header + values for column B
$ws.Range("B1").Value = "expected Outcome  "
$ws.Range("B2").Value = "popup error message containing SyntaxError "
$ws.Range("B3").Value = "the user is able to see the output in the console"

# Widen the new column to fit its content
$ws.Columns("B").ColumnWidth = 39

# Match the resulting active-cell selection on this sheet
$ws.Range("D4").Select()
